# The upstream PanelApp query was re-run, producing fresh per-row query
# timestamps in the "data" sheet's time_taken column (F2:F38).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$newTimes = @(
  "2021-10-05 14:23:06.350618",
  "2021-10-05 14:23:06.350627",
  "2021-10-05 14:23:06.350631",
  "2021-10-05 14:23:06.350633",
  "2021-10-05 14:23:06.350636",
  "2021-10-05 14:23:06.350639",
  "2021-10-05 14:23:06.350642",
  "2021-10-05 14:23:06.350644",
  "2021-10-05 14:23:06.350647",
  "2021-10-05 14:23:06.350650",
  "2021-10-05 14:23:06.350653",
  "2021-10-05 14:23:06.350655",
  "2021-10-05 14:23:06.350658",
  "2021-10-05 14:23:06.350661",
  "2021-10-05 14:23:06.350663",
  "2021-10-05 14:23:06.350666",
  "2021-10-05 14:23:06.350669",
  "2021-10-05 14:23:06.350672",
  "2021-10-05 14:23:06.350675",
  "2021-10-05 14:23:06.350678",
  "2021-10-05 14:23:06.350681",
  "2021-10-05 14:23:06.350684",
  "2021-10-05 14:23:06.350687",
  "2021-10-05 14:23:06.350689",
  "2021-10-05 14:23:06.350693",
  "2021-10-05 14:23:06.350695",
  "2021-10-05 14:23:06.350698",
  "2021-10-05 14:23:06.350701",
  "2021-10-05 14:23:06.350704",
  "2021-10-05 14:23:06.350706",
  "2021-10-05 14:23:06.350709",
  "2021-10-05 14:23:06.350712",
  "2021-10-05 14:23:06.350715",
  "2021-10-05 14:23:06.350717",
  "2021-10-05 14:23:06.350720",
  "2021-10-05 14:23:06.350723",
  "2021-10-05 14:23:06.350726"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Split the panel-level metadata (name/id/version/request info) out of the
# single "data" tab into its own "metadata" tab, placed right after "data".
$after = $wb.Worksheets.Item("data")
$meta = $wb.Worksheets.Add($null, $after)
$meta.Name = "metadata"

# Reuse the bold/bordered/centred header style from the "data" sheet for
# the header row and the row-index column, instead of re-deriving a
# lookalike style (which would mint a brand-new style entry).
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Vascular skin disorders"
$meta.Cells.Item(2, 3).Value = 563
# "1.47" is a version label, not a number -- force text storage so it
# round-trips as a string, then drop back to the default (unstyled) look.
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.47"
$meta.Cells.Item(2, 4).Style = "Normal"
$meta.Cells.Item(2, 5).Value = "2021-03-22T15:24:57.558379Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:23:06.347330"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/563/?format=json"

[void]$meta.Range("A1").Select()

# Leave the workbook's active tab on "data" (unchanged from before).
$ws.Activate()
[void]$ws.Range("A1").Select()
